$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that looks numeric (e.g. "35.29").
# Force text storage (matching the source inlineStr type) so Excel
# does not silently coerce them to floating-point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.970.60"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.085.84"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.86"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.39"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.13"
$ws.Range("E7").Value = "  +8.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.367"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.084.84"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.748"
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.29"
$ws.Range("E14").Value = "  -3.63%  "
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.47"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.046.89"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.662.39"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.088.93"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.68"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.58"
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.81"
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "443.63"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.09"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "91.13"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.61"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.87"
$ws.Range("E27").Value = "  -5.85%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.186"
$ws.Range("E30").Value = "  +16.37%  "
$ws.Range("E31").Value = "  +26.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.24"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").Value = "  +13.76%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.169"
$ws.Range("E34").Value = "  +14.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.112"
$ws.Range("E35").Value = "  +32.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.72"
$ws.Range("E36").Value = "  +7.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.46"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.18"
$ws.Range("E38").Value = "  +27.82%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "493.91"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  -5.38%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.418"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.11"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.98"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.691"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.47"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.33"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.17"
$ws.Range("E51").Value = "  -2.23%  "
